$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), copying the style from H1 (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I and J columns for rows 2-71
$data = @(
    @(7,7),
    @(8,8),
    @(5,6),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(6,7),
    @(6,6),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(7,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(6,6),
    @(7,8),
    @(9,9),
    @(10,11),
    @(9,9),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,6),
    @(9,9),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(5,5),
    @(4,4),
    @(1,1),
    @(4,4)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}

Write-Output "done"
